# Refresh current Universalis price snapshots + recomputed Leve profit
# columns (H-N) for the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-leve
# sheets, per the scheduled market-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$updates = @{
  "H17" = 3208431.2
  "J17" = 3208431.2
  "L17" = 9625293.600000001
  "N17" = -9625629.600000001
  "H33" = 6062058.5
  "I33" = 1556.6428
  "J33" = 90909090
  "K33" = 1556.6428
  "L33" = 90909090
  "M33" = -1327.6428
  "N33" = -90909548
  "H51" = 2644.4443
  "I51" = 2900
  "J51" = 2571.4285
  "K51" = 2900
  "L51" = 2571.4285
  "M51" = -2416
  "N51" = -3539.4285
  "H106" = 444447780
  "I106" = 166671680
  "K106" = 166671680
  "M106" = -166671049
  "H135" = 2004.8235
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")

$updates = @{
  "H104" = 38000
  "J104" = 38000
  "L104" = 38000
  "N104" = -44988
  "H129" = 49999
  "J129" = 49999
  "L129" = 49999
  "N129" = -59999
  "H132" = 6669.033
  "I132" = 1489.2142
  "J132" = 11201.375
  "K132" = 4467.642599999999
  "L132" = 33604.125
  "M132" = -1937.642599999999
  "N132" = -38664.125
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("BSM")

$updates = @{
  "H64" = 237.18182
  "I64" = 183.66667
  "J64" = 301.4
  "K64" = 183.66667
  "L64" = 301.4
  "M64" = 41.33332999999999
  "N64" = -751.4
  "H67" = 237.18182
  "I67" = 183.66667
  "J67" = 301.4
  "K67" = 183.66667
  "L67" = 301.4
  "M67" = 596.3333299999999
  "N67" = -1861.4
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")

$updates = @{
  "H58" = 863.5
  "J58" = 0
  "L58" = 0
  "H132" = 2973.6428
  "I132" = 2423.4
  "J132" = 4349.25
  "K132" = 7270.200000000001
  "L132" = 13047.75
  "M132" = -4740.200000000001
  "N132" = -18107.75
  "H134" = 1589.1
  "I134" = 1543.4445
  "K134" = 4630.333500000001
  "M134" = -2095.333500000001
  "H136" = 863.5
  "J136" = 0
  "L136" = 0
  "H138" = 110000
  "J138" = 110000
  "L138" = 110000
  "N138" = -120280
  "H139" = 0
  "J139" = 0
  "L139" = 0
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}
$clears = @("N58", "N136", "N139")
foreach ($addr in $clears) {
  $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("CUL")

$updates = @{
  "H34" = 1963.48
  "I34" = 285.375
  "J34" = 2753.1765
  "K34" = 856.125
  "L34" = 8259.529500000001
  "M34" = -772.125
  "N34" = -8427.529500000001
  "H39" = 2068
  "J39" = 2068
  "L39" = 6204
  "N39" = -6792
  "H55" = 4129.8887
  "I55" = 1400
  "J55" = 4471.125
  "K55" = 4200
  "L55" = 13413.375
  "M55" = -4023
  "N55" = -13767.375
  "H82" = 3664.7058
  "H85" = 3664.7058
  "H94" = 3885.7144
  "J94" = 3885.7144
  "L94" = 11657.1432
  "N94" = -13009.1432
  "H97" = 33333834
  "I97" = 50000250
  "J97" = 1000
  "K97" = 150000750
  "L97" = 3000
  "M97" = -150000254
  "N97" = -3992
  "H103" = 1173.5
  "I103" = 598
  "J103" = 2900
  "K103" = 1794
  "L103" = 8700
  "M103" = -915
  "N103" = -10458
  "H109" = 1592.5
  "J109" = 3300
  "L109" = 9900
  "N109" = -11980
  "H118" = 1197.4
  "I118" = 746.75
  "J118" = 3000
  "K118" = 2240.25
  "L118" = 9000
  "M118" = -997.25
  "N118" = -11486
  "H121" = 1001.79486
  "I121" = 433.33334
  "J121" = 1049.1666
  "K121" = 1300.00002
  "L121" = 3147.4998
  "M121" = 9.99998000000005
  "N121" = -5767.4998
  "H129" = 15152922
  "J129" = 2287.5
  "L129" = 6862.5
  "N129" = -16862.5
  "H131" = 977.4545000000001
  "I131" = 0
  "J131" = 977.4545000000001
  "K131" = 0
  "L131" = 2932.3635
  "N131" = -13012.3635
  "H137" = 14428.044
  "J137" = 13526.625
  "L137" = 40579.875
  "N137" = -50779.875
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}
$clears = @("M131")
foreach ($addr in $clears) {
  $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("GSM")

$updates = @{
  "H122" = 2402175.2
  "I122" = 3602249
  "K122" = 10806747
  "M122" = -10804297
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")

$updates = @{
  "H40" = 200004800
  "I40" = 250003740
  "K40" = 250003740
  "M40" = -250003604
  "H93" = 31263818
  "I93" = 19181.818
  "J93" = 100002020
  "K93" = 19181.818
  "L93" = 100002020
  "M93" = -17933.818
  "N93" = -100004516
  "H106" = 70000
  "J106" = 70000
  "L106" = 70000
  "N106" = -72524
  "H132" = 14450328
  "I132" = 19703852
  "J132" = 3137.875
  "K132" = 59111556
  "L132" = 9413.625
  "M132" = -59109026
  "N132" = -14473.625
  "H136" = 7434.5483
  "I136" = 6350.778
  "J136" = 14750
  "K136" = 19052.334
  "L136" = 44250
  "M136" = -16502.334
  "N136" = -49350
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("WVR")

$updates = @{
  "H64" = 25114
  "J64" = 25114
  "L64" = 25114
  "N64" = -25610
  "H67" = 25114
  "J67" = 25114
  "L67" = 25114
  "N67" = -26830
  "H104" = 41685
  "J104" = 41685
  "L104" = 41685
  "N104" = -48673
  "H105" = 48896
  "J105" = 48896
  "L105" = 48896
  "N105" = -55884
  "H132" = 1385.7587
  "I132" = 866.8
  "K132" = 2600.4
  "M132" = -70.39999999999964
}
foreach ($addr in $updates.Keys) {
  $ws.Range($addr).Value = $updates[$addr]
}

